{"js": "// Update the date line and the table of division problems.\n// Both the date paragraph and each table cell are addressed positionally\n// (document order) so duplicate values (e.g. the two \"88\u00f77=\" cells) are\n// handled unambiguously.\n\nconst body = context.document.body;\n\n// --- 1. Update the date/day-of-week line at the top of the document ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.load(\"text\");\nawait context.sync();\n\nif (firstParagraph.text.trim() === \"2024-10-17 Thursday\") {\n  firstParagraph.insertText(\"2024-10-18 Friday\", Word.InsertLocation.replace);\n}\n\n// --- 2. Update the worksheet table, cell by cell, in document order ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Map of (row index -> [old -> new, ...]) values, in left-to-right order,\n// for each of the five populated rows (0, 4, 8, 12, 16).\nconst rowUpdates = {\n  0: [\"94\u00f73=\", \"48\u00f79=\", \"62\u00f79=\", \"68\u00f77=\", \"18\u00f79=\"],\n  4: [\"21\u00f75=\", \"84\u00f77=\", \"68\u00f73=\", \"66\u00f73=\", \"36\u00f74=\"],\n  8: [\"34\u00f75=\", \"65\u00f74=\", \"21\u00f74=\", \"12\u00f78=\", \"48\u00f73=\"],\n  12: [\"70\u00f74=\", \"89\u00f79=\", \"13\u00f76=\", \"33\u00f76=\", \"54\u00f76=\"],\n  16: [\"50\u00f76=\", \"33\u00f72=\", \"70\u00f76=\", \"84\u00f75=\", \"25\u00f73=\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newValues = rowUpdates[rowIndex];\n  const row = rows.items[rowIndex];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let c = 0; c < cells.items.length; c++) {\n    const cellBody = cells.items[c].body;\n    const cellParagraphs = cellBody.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n\n    // Replace the text of the existing paragraph's range in place so the\n    // run-level formatting (font/size) and paragraph formatting (jc) are\n    // preserved, instead of wiping the cell body and inserting a fresh run.\n    const cellParagraph = cellParagraphs.items[0];\n    const cellRange = cellParagraph.getRange();\n    cellRange.insertText(newValues[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the table of division problems.\n# Both the date paragraph and each table cell are addressed positionally\n# (row/column index, 1-based like Word's COM model) so duplicate values\n# (e.g. the two \"88\u00f77=\" cells) are handled unambiguously.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date/day-of-week line at the top of the document ---\n$firstParagraph = $d.Paragraphs.Item(1)\nif ($firstParagraph.Range.Text.Trim() -eq \"2024-10-17 Thursday\") {\n    $firstParagraph.Range.Text = \"2024-10-18 Friday\"\n}\n\n# --- 2. Update the worksheet table, cell by cell, in document order ---\n$table = $d.Tables.Item(1)\n\n# Map of (1-based row index -> new values for the 5 columns, left to right)\n# for each of the five populated rows (1, 5, 9, 13, 17).\n$rowUpdates = @{\n    1  = @(\"94\u00f73=\", \"48\u00f79=\", \"62\u00f79=\", \"68\u00f77=\", \"18\u00f79=\")\n    5  = @(\"21\u00f75=\", \"84\u00f77=\", \"68\u00f73=\", \"66\u00f73=\", \"36\u00f74=\")\n    9  = @(\"34\u00f75=\", \"65\u00f74=\", \"21\u00f74=\", \"12\u00f78=\", \"48\u00f73=\")\n    13 = @(\"70\u00f74=\", \"89\u00f79=\", \"13\u00f76=\", \"33\u00f76=\", \"54\u00f76=\")\n    17 = @(\"50\u00f76=\", \"33\u00f72=\", \"70\u00f76=\", \"84\u00f75=\", \"25\u00f73=\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $newValues = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $newValues.Length; $col++) {\n        $cell = $table.Cell($rowIndex, $col)\n        $cell.Range.Text = $newValues[$col - 1]\n    }\n}\n"}
